$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(95878, "Emanuella Albuquerque", "Recursos Humanos", "Viagem de negocios", 4, "2023-06-06", 4010.87),
    @(65745, "Gael Andrade", "Juridico", "Outros", 6, "2023-06-06", 5138.41),
    @(18925, "Srta. Elisa da Mata", "Marketing", "Consulta medica", 2, "2023-06-13", 7303.7),
    @(32183, "Sr. Benicio Ramos", "Marketing", "Outros", 7, "2023-06-03", 6619.7),
    @(11573, "Fernando Cassiano", "TI", "Doenca", 1, "2023-06-07", 2082.89),
    @(24073, "Dr. Anthony Freitas", "Recursos Humanos", "Outros", 2, "2023-06-11", 6944.59),
    @(76366, "Isabela Sales", "Vendas", "Doenca", 2, "2023-06-24", 5992.03),
    @(78776, "Carolina Cassiano", "Marketing", "Outros", 4, "2023-06-03", 3167.26),
    @(60031, "Luiza Carvalho", "P&D", "Consulta medica", 2, "2023-06-01", 9220.34),
    @(91965, "Vitória Albuquerque", "Financeiro", "Doenca", 2, "2023-06-03", 8193.469999999999)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $row++
}
